# Forest data update - 2026-01-06 12:20
# -----------------------------------------------------------------------
# 5 listings that were sitting in the "New" sheet (rows 2-6) have now been
# vetted, so they move down to the end of "Previously added" (rows 355-359,
# unchanged). The "New" sheet is then repopulated with 4 freshly scraped
# listings (rows 2-5).
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

$xlPasteFormats = -4122

# -------------------------------------------------------------------
# 1) The 5 rows currently on "New" (rows 2-6) - these are about to be
#    overwritten on that sheet, so capture their values/links first.
# -------------------------------------------------------------------
$movingRows = @()
for ($r = 2; $r -le 6; $r++) {
    $movingRows += , @(
        $wsNew.Cells.Item($r, 1).Value(),
        $wsNew.Cells.Item($r, 2).Value(),
        $wsNew.Cells.Item($r, 3).Value(),
        $wsNew.Cells.Item($r, 4).Value(),
        $wsNew.Cells.Item($r, 5).Value(),
        $wsNew.Cells.Item($r, 6).Value()
    )
}

# -------------------------------------------------------------------
# 2) Append those 5 rows to the bottom of "Previously added"
#    (rows 355-359), matching the formatting of the last existing row.
# -------------------------------------------------------------------
$lastRow = 354
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $movingRows.Count

$wsPrev.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$wsPrev.Range("A" + $firstNewRow + ":F" + $lastNewRow).PasteSpecial($xlPasteFormats)

$r = $firstNewRow
foreach ($row in $movingRows) {
    # cadastre numbers (col E) must stay text, not be coerced to numbers
    $wsPrev.Cells.Item($r, 5).NumberFormat = "@"

    $wsPrev.Cells.Item($r, 1).Value = $row[0]
    $wsPrev.Cells.Item($r, 2).Value = $row[1]
    $wsPrev.Cells.Item($r, 3).Value = $row[2]
    $wsPrev.Cells.Item($r, 4).Value = $row[3]
    $wsPrev.Cells.Item($r, 5).Value = $row[4]
    $wsPrev.Cells.Item($r, 6).Value = $row[5]
    $wsPrev.Hyperlinks.Add($wsPrev.Cells.Item($r, 1), $row[0])
    $r = $r + 1
}

# restore the usual hyperlink-cell styling (Hyperlinks.Add applies the
# built-in "Hyperlink" style; put back the sheet's own custom style) and
# undo the text NumberFormat override on column E
$wsPrev.Range("A" + $lastRow).Copy()
$wsPrev.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial($xlPasteFormats)
$wsPrev.Range("E" + $lastRow).Copy()
$wsPrev.Range("E" + $firstNewRow + ":E" + $lastNewRow).PasteSpecial($xlPasteFormats)

# -------------------------------------------------------------------
# 3) Replace the content of "New": drop all 5 old hyperlinks/rows and
#    write the 4 newly scraped listings into rows 2-5.
# -------------------------------------------------------------------
$wsNew.Hyperlinks.Delete()
$wsNew.Rows.Item(6).Delete()

$newListings = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/priekulu-pag/bdjkcx.html", "36 000 €", "Cēsis un raj.", "11 ha.", "42720030144", 46028.57847222222),
    @("https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/liepas-pag/lkdcm.html", "120 000 €", "Cēsis un raj.", "22 ha.", "42600030079", 46027.75486111111),
    @("https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/galgauskas-pag/khelk.html", "29 600 €", "Gulbene un raj.", "5 ha.", "50560060064", 46028.57847222222),
    @("https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/liepupes-pag/hoxnl.html", "58 000 €", "Limbaži un raj.", "4.60 ha.", "", 46028.51527777778)
)

$r = 2
foreach ($row in $newListings) {
    # columns that look numeric (cadastre numbers) must stay text
    $wsNew.Cells.Item($r, 5).NumberFormat = "@"

    $wsNew.Cells.Item($r, 1).Value = $row[0]
    $wsNew.Cells.Item($r, 2).Value = $row[1]
    $wsNew.Cells.Item($r, 3).Value = $row[2]
    $wsNew.Cells.Item($r, 4).Value = $row[3]
    $wsNew.Cells.Item($r, 5).Value = $row[4]
    $wsNew.Cells.Item($r, 6).Value = $row[5]
    $wsNew.Hyperlinks.Add($wsNew.Cells.Item($r, 1), $row[0])
    $r = $r + 1
}

# restore original per-column styling that Hyperlinks.Add / NumberFormat
# touched (col A hyperlink style, col E plain text style)
$wsPrev.Range("A" + $lastRow).Copy()
$wsNew.Range("A2:A5").PasteSpecial($xlPasteFormats)
$wsNew.Range("B2").Copy()
$wsNew.Range("E2:E5").PasteSpecial($xlPasteFormats)

Write-Host "Forest data updated"
